# Weekly refresh of Fruta/Hortaliza data: each price record (columns D..T)
# is re-shuffled across the existing rows (2..29) to reflect the updated
# daily-logic snapshot. Columns A..C and E..J are constant for every row
# in this sheet and are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of destination row -> source row (i.e. new row R gets the old
# content that used to live in row Map[R]).
$map = @{
    2  = 19
    3  = 29
    4  = 16
    5  = 20
    6  = 18
    7  = 22
    8  = 27
    9  = 28
    10 = 15
    11 = 24
    12 = 26
    13 = 25
    14 = 23
    15 = 7
    16 = 3
    17 = 4
    18 = 13
    19 = 12
    20 = 2
    21 = 6
    22 = 14
    23 = 5
    24 = 9
    25 = 11
    26 = 8
    27 = 21
    28 = 17
    29 = 10
}

$firstRow = 2
$lastRow = 29
$firstCol = 4   # D
$lastCol = 20   # T

# 1) Snapshot the current D:T content of every row before overwriting
#    anything (needed because the permutation reuses rows as sources).
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @{}
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

# 2) Write the permuted values back into the sheet.
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $srcRow = $map[$r]
    $srcVals = $snapshot[$srcRow]
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $ws.Cells.Item($r, $c).Value = $srcVals[$c]
    }
}
